# Auto-generated edit script applying scheduled-runner price/profit updates
# to the Seraph_Profits sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 573.7091
$ws.Range("I15").Value = 573.7091
$ws.Range("K15").Value = 1721.1273
$ws.Range("M15").Value = -1552.1273
$ws.Range("H32").Value = 933
$ws.Range("J32").Value = 933
$ws.Range("L32").Value = 933
$ws.Range("N32").Value = -1585
$ws.Range("H41").Value = 395.9091
$ws.Range("I41").Value = 531.4
$ws.Range("J41").Value = 283
$ws.Range("K41").Value = 531.4
$ws.Range("L41").Value = 283
$ws.Range("M41").Value = -91.39999999999998
$ws.Range("N41").Value = -1163
$ws.Range("H58").Value = 3613.8
$ws.Range("J58").Value = 5999.6665
$ws.Range("L58").Value = 17998.9995
$ws.Range("N58").Value = -18298.9995
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()
$ws.Range("H112").Value = 2974.111
$ws.Range("I112").Value = 1420
$ws.Range("J112").Value = 3065.5293
$ws.Range("K112").Value = 4260
$ws.Range("L112").Value = 9196.5879
$ws.Range("M112").Value = -3152
$ws.Range("N112").Value = -11412.5879
$ws.Range("H116").Value = 5555
$ws.Range("I116").Value = 3260
$ws.Range("K116").Value = 3260
$ws.Range("M116").Value = 182
$ws.Range("H137").Value = 4565.2173
$ws.Range("I137").Value = 1984.2307
$ws.Range("K137").Value = 5952.6921
$ws.Range("M137").Value = -3402.6921
$ws.Range("H138").Value = 6184.0586
$ws.Range("I138").Value = 6534.579
$ws.Range("J138").Value = 5740.067
$ws.Range("K138").Value = 19603.737
$ws.Range("L138").Value = 17220.201
$ws.Range("M138").Value = -14463.737
$ws.Range("N138").Value = -27500.201

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18052.475
$ws.Range("I32").Value = 8625.6
$ws.Range("K32").Value = 8625.6
$ws.Range("M32").Value = -8338.6
$ws.Range("H74").Value = 4960.533
$ws.Range("I74").Value = 1666.3334
$ws.Range("K74").Value = 1666.3334
$ws.Range("M74").Value = -792.3334
$ws.Range("H77").Value = 4960.533
$ws.Range("I77").Value = 1666.3334
$ws.Range("K77").Value = 8331.666999999999
$ws.Range("M77").Value = -3963.666999999999
$ws.Range("H102").Value = 1078.3
$ws.Range("I102").Value = 1039.7142
$ws.Range("K102").Value = 1039.7142
$ws.Range("M102").Value = 582.2858000000001
$ws.Range("H110").Value = 5244.8
$ws.Range("I110").Value = 5282.4614
$ws.Range("K110").Value = 5282.4614
$ws.Range("M110").Value = -3237.4614
$ws.Range("H122").Value = 3464.2856
$ws.Range("I122").Value = 2648.75
$ws.Range("J122").Value = 4551.6665
$ws.Range("K122").Value = 7946.25
$ws.Range("L122").Value = 13654.9995
$ws.Range("M122").Value = -5496.25
$ws.Range("N122").Value = -18554.9995
$ws.Range("H132").Value = 3880.681
$ws.Range("I132").Value = 1540.9117
$ws.Range("K132").Value = 4622.7351
$ws.Range("M132").Value = -2092.7351

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 682.1111
$ws.Range("I22").Value = 682.1111
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 682.1111
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -509.1111
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 1201
$ws.Range("I94").Value = 1351.5
$ws.Range("K94").Value = 1351.5
$ws.Range("M94").Value = -900.5
$ws.Range("H105").Value = 4410.9644
$ws.Range("I105").Value = 3753.3125
$ws.Range("K105").Value = 3753.3125
$ws.Range("M105").Value = -2006.3125
$ws.Range("H134").Value = 2542.111
$ws.Range("I134").Value = 1856.591
$ws.Range("J134").Value = 5558.4
$ws.Range("K134").Value = 5569.772999999999
$ws.Range("L134").Value = 16675.2
$ws.Range("M134").Value = -3034.772999999999
$ws.Range("N134").Value = -21745.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 750
$ws.Range("I16").Value = 750
$ws.Range("K16").Value = 750
$ws.Range("M16").Value = -463
$ws.Range("H22").Value = 279.6
$ws.Range("I22").Value = 132.66667
$ws.Range("J22").Value = 500
$ws.Range("K22").Value = 132.66667
$ws.Range("L22").Value = 500
$ws.Range("M22").Value = 217.33333
$ws.Range("N22").Value = -1200
$ws.Range("H58").Value = 3749.8096
$ws.Range("I58").Value = 1536.3572
$ws.Range("K58").Value = 1536.3572
$ws.Range("M58").Value = -1333.3572
$ws.Range("H94").Value = 1572.3334
$ws.Range("J94").Value = 1572.3334
$ws.Range("L94").Value = 1572.3334
$ws.Range("N94").Value = -2474.3334
$ws.Range("H109").Value = 16353.637
$ws.Range("J109").Value = 16353.637
$ws.Range("L109").Value = 16353.637
$ws.Range("N109").Value = -18433.637
$ws.Range("H113").Value = 750
$ws.Range("I113").Value = 750
$ws.Range("K113").Value = 750
$ws.Range("M113").Value = 1420
$ws.Range("H122").Value = 916.4
$ws.Range("I122").Value = 895.5
$ws.Range("K122").Value = 2686.5
$ws.Range("M122").Value = -236.5
$ws.Range("H136").Value = 3749.8096
$ws.Range("I136").Value = 1536.3572
$ws.Range("K136").Value = 4609.071599999999
$ws.Range("M136").Value = -2059.071599999999
$ws.Range("H141").Value = 92490
$ws.Range("J141").Value = 92490
$ws.Range("L141").Value = 92490
$ws.Range("N141").Value = -102850

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1263.4
$ws.Range("I131").Value = 799.6667
$ws.Range("K131").Value = 2399.0001
$ws.Range("M131").Value = 2640.9999
$ws.Range("H132").Value = 7866.6665
$ws.Range("J132").Value = 7866.6665
$ws.Range("L132").Value = 70799.9985
$ws.Range("N132").Value = -75859.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 4657.857
$ws.Range("I113").Value = 3439.4
$ws.Range("K113").Value = 3439.4
$ws.Range("M113").Value = -1269.4

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1473.3334
$ws.Range("I22").Value = 792.2
$ws.Range("J22").Value = 2324.75
$ws.Range("K22").Value = 792.2
$ws.Range("L22").Value = 2324.75
$ws.Range("M22").Value = -497.2
$ws.Range("N22").Value = -2914.75
$ws.Range("H27").Value = 1473.3334
$ws.Range("I27").Value = 792.2
$ws.Range("J27").Value = 2324.75
$ws.Range("K27").Value = 792.2
$ws.Range("L27").Value = 2324.75
$ws.Range("M27").Value = -685.2
$ws.Range("N27").Value = -2538.75
$ws.Range("H33").Value = 238336670
$ws.Range("I33").Value = 238336670
$ws.Range("K33").Value = 238336670
$ws.Range("M33").Value = -238336380
$ws.Range("H46").Value = 3449.875
$ws.Range("I46").Value = 1866.6666
$ws.Range("J46").Value = 4399.8
$ws.Range("K46").Value = 1866.6666
$ws.Range("L46").Value = 4399.8
$ws.Range("M46").Value = -1678.6666
$ws.Range("N46").Value = -4775.8
$ws.Range("H93").Value = 3000
$ws.Range("J93").Value = 3000
$ws.Range("L93").Value = 3000
$ws.Range("N93").Value = -5496
$ws.Range("H132").Value = 4762.8125
$ws.Range("I132").Value = 3411.7778
$ws.Range("K132").Value = 10235.3334
$ws.Range("M132").Value = -7705.3334
$ws.Range("H136").Value = 3499.5
$ws.Range("J136").Value = 4999
$ws.Range("L136").Value = 14997
$ws.Range("N136").Value = -20097

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 19999.75
$ws.Range("J54").Value = 19999.75
$ws.Range("L54").Value = 19999.75
$ws.Range("N54").Value = -21039.75
$ws.Range("H62").Value = 7681.1816
$ws.Range("J62").Value = 8055.5557
$ws.Range("L62").Value = 8055.5557
$ws.Range("N62").Value = -9303.555700000001
$ws.Range("H65").Value = 7681.1816
$ws.Range("J65").Value = 8055.5557
$ws.Range("L65").Value = 40277.7785
$ws.Range("N65").Value = -46517.7785
$ws.Range("H113").Value = 1555.4286
$ws.Range("J113").Value = 2679.6
$ws.Range("L113").Value = 8038.799999999999
$ws.Range("N113").Value = -12378.8
$ws.Range("H132").Value = 2598.4
$ws.Range("I132").Value = 2069.5
$ws.Range("J132").Value = 3832.5
$ws.Range("K132").Value = 6208.5
$ws.Range("L132").Value = 11497.5
$ws.Range("M132").Value = -3678.5
$ws.Range("N132").Value = -16557.5
$ws.Range("H136").Value = 61861.766
$ws.Range("I136").Value = 1705.6364
$ws.Range("K136").Value = 5116.9092
$ws.Range("M136").Value = -2566.9092
